$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("A1").Value = "Segment"
$ws.Range("B1").Value = "Country"
$ws.Range("C1").Value = "Product"
$ws.Range("D1").Value = "Discount Band"
$ws.Range("E1").Value = "Units Sold"
$ws.Range("F1").Value = "Manufacturing P"

# Data rows (Financial Sample style data)
$ws.Range("A2").Value = "Government"
$ws.Range("B2").Value = "Canada"
$ws.Range("C2").Value = "Carretera"
$ws.Range("D2").Value = "None"
$ws.Range("E2").Value = 1618.5
$ws.Range("F2").Value = "US$  3.00"

$ws.Range("A3").Value = "Government"
$ws.Range("B3").Value = "Germany"
$ws.Range("C3").Value = "Carretera"
$ws.Range("D3").Value = "None"
$ws.Range("E3").Value = 1321
$ws.Range("F3").Value = "US$  3.00"

$ws.Range("A4").Value = "Midmarket"
$ws.Range("B4").Value = "France"
$ws.Range("C4").Value = "Carretera"
$ws.Range("D4").Value = "None"
$ws.Range("E4").Value = 2178
$ws.Range("F4").Value = "US$  3.00"

$ws.Range("A5").Value = "Midmarket"
$ws.Range("B5").Value = "Germany"
$ws.Range("C5").Value = "Carretera"
$ws.Range("D5").Value = "None"
$ws.Range("E5").Value = 888
$ws.Range("F5").Value = "US$  3.00"

$ws.Range("A6").Value = "Midmarket"
$ws.Range("B6").Value = "Mexico"
$ws.Range("C6").Value = "Carretera"
$ws.Range("D6").Value = "None"
$ws.Range("E6").Value = 2470
$ws.Range("F6").Value = "US$  3.00"

$ws.Range("A7").Value = "Government"
$ws.Range("B7").Value = "Germany"
$ws.Range("C7").Value = "Carretera"
$ws.Range("D7").Value = "None"
$ws.Range("E7").Value = 1513
$ws.Range("F7").Value = "US$  3.00"

$ws.Range("A8").Value = "Midmarket"
$ws.Range("B8").Value = "Germany"
$ws.Range("C8").Value = "Montana"
$ws.Range("D8").Value = "None"
$ws.Range("E8").Value = 921
$ws.Range("F8").Value = "US$  3.00"

$ws.Range("A9").Value = "Channel Partners"
$ws.Range("B9").Value = "Canada"
$ws.Range("C9").Value = "Montana"
$ws.Range("D9").Value = "None"
$ws.Range("E9").Value = 2518
$ws.Range("F9").Value = "US$  3.00"

# Column widths (characters) to match the saved layout
$ws.Range("A1:B1").EntireColumn.ColumnWidth = 17.333333333333332
$ws.Range("D1").EntireColumn.ColumnWidth = 14
$ws.Range("E1").EntireColumn.ColumnWidth = 9.166666666666666
$ws.Range("F1").EntireColumn.ColumnWidth = 16

# Selection matching the saved state
$ws.Range("G9").Select() | Out-Null
